# Refresh the cryptos list (Price / Volume(1h) columns) with latest scraped
# values, as produced by the "Updated cryptos list ... with GitHub Actions"
# job. Price/volume cells are stored as plain text in this sheet, so any
# value that would otherwise be auto-parsed as a number is written with a
# leading "'" text-qualifier to keep it as text (this only sets Excel's
# "quote prefix" flag - it is not part of the stored string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.027.44"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.471.86"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D5").Value = "'584.34"
$ws.Range("D6").Value = "'171.50"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "2.470.93"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +3.82%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").Value = "'4.92"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "2.918.80"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "'25.35"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "66.931.58"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "2.479.81"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'10.92"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "'7.48"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'349.69"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "'3.98"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'68.36"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'1.78"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "0.0₃0898"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "'506.77"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "'7.67"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'161.91"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "'18.12"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.328"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'4.80"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "'2.36"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").Value = "'142.64"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'3.47"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'0.512"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "0.0₆0251"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "'0.0732"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'1.56"
$ws.Range("E51").Value = "  -1.08%  "
